$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.490.82"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'2.134.75"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'351.66"
$ws.Range("E5").Value = "  +5.18%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.5247"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "'0.4538"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "'53.59"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'0.09137"
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("D11").Value = "'1.189"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'25.40"
$ws.Range("E12").Value = "  +5.31%  "
$ws.Range("D13").Value = "'2.145.29"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").Value = "'6.873"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "'8.150"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "'101.66"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "'0.00001167"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'0.06716"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").Value = "'20.36"
$ws.Range("E20").Value = "  +6.16%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'6.366"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").Value = "'30.579.90"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'12.85"
$ws.Range("E24").Value = "  +4.05%  "
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "'2.391.55"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").Value = "'22.48"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "'2.613"
$ws.Range("E28").Value = "  +3.77%  "
$ws.Range("D29").Value = "'164.61"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").Value = "'135.70"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "'1.223"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Value = "'1.715"
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("D33").Value = "'0.1081"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").Value = "'6.380"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'4.025"
$ws.Range("E35").Value = "  +2.35%  "
$ws.Range("D36").Value = "'6.114"
$ws.Range("E36").Value = "  +4.11%  "
$ws.Range("D37").Value = "'10.43"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'0.02641"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("D39").Value = "'0.06965"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'0.2341"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "'12.68"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'0.6968"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "'1.273"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").Value = "'14.84"
$ws.Range("E44").Value = "  +5.93%  "

# Rows 45 and 46 swap coin identity (NEARProtocol <-> Decentraland) with updated data
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6505"
$ws.Range("E45").Value = "  +2.41%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.346"
$ws.Range("E46").Value = "  +1.01%  "

$ws.Range("D47").Value = "'0.00000000373"
$ws.Range("E47").Value = "  +9.48%  "
$ws.Range("D48").Value = "'3.739"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("D49").Value = "'1.249"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'83.77"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'0.07289"
$ws.Range("E51").Value = "  +2.28%  "
